$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "Rgmb", "Neo1", "ECs", 2, 0.6666666666666666, 9.249291999999999, 27.747876, 0.1285586120341026, 0.1285586120341026, 3, 1, 2.373057666666667, 7.119173, 0.03663411841135517, 0.03663411841135517, 21.94910329183866, 197.541929626548, 0.004709631416056785, 0.004709631416056783),
  @("ECs", "Rgmb", "Neo1", "FAPs", 2, 0.6666666666666666, 9.249291999999999, 27.747876, 0.1285586120341026, 0.1285586120341026, 3, 1, 43.26393666666667, 129.79181, 0.6678877639810287, 0.6678877639810287, 400.1607832995066, 3601.44704969556, 0.08586272393196137, 0.08586272393196134),
  @("ECs", "Rgmb", "Neo1", "M2", 2, 0.6666666666666666, 9.249291999999999, 27.747876, 0.1285586120341026, 0.1285586120341026, 3, 1, 0.4664283333333333, 1.399285, 0.007200495391983466, 0.007200495391983466, 4.314131852073332, 38.82718666865999, 0.0009256856935513461, 0.0009256856935513459),
  @("ECs", "Rgmb", "Neo1", "sCs", 2, 0.6666666666666666, 9.249291999999999, 27.747876, 0.1285586120341026, 0.1285586120341026, 3, 1, 18.67383333333333, 56.0215, 0.2882776222156329, 0.2882776222156329, 172.7197372593333, 1554.477635334, 0.03706057099253315, 0.03706057099253314),
  @("FAPs", "Rgmb", "Neo1", "ECs", 3, 1, 35.42317866666667, 106.269536, 0.492357110492641, 0.492357110492641, 3, 1, 2.373057666666667, 7.119173, 0.03663411841135517, 0.03663411841135517, 84.06124571263643, 756.551211413728, 0.01803706868646009, 0.01803706868646009),
  @("FAPs", "Rgmb", "Neo1", "FAPs", 3, 1, 35.42317866666667, 106.269536, 0.492357110492641, 0.492357110492641, 3, 1, 43.26393666666667, 129.79181, 0.6678877639810287, 0.6678877639810287, 1532.546158366684, 13792.91542530016, 0.3288392896070903, 0.3288392896070902),
  @("FAPs", "Rgmb", "Neo1", "M2", 3, 1, 35.42317866666667, 106.269536, 0.492357110492641, 0.492357110492641, 3, 1, 0.4664283333333333, 1.399285, 0.007200495391983466, 0.007200495391983466, 16.52237418686222, 148.70136768176, 0.003545215105312556, 0.003545215105312556),
  @("FAPs", "Rgmb", "Neo1", "sCs", 3, 1, 35.42317866666667, 106.269536, 0.492357110492641, 0.492357110492641, 3, 1, 18.67383333333333, 56.0215, 0.2882776222156329, 0.2882776222156329, 661.4865345582223, 5953.378811024, 0.1419355370937782, 0.1419355370937782),
  @("M2", "Rgmb", "Neo1", "ECs", 3, 1, 3.703826, 11.111478, 0.0514805597850973, 0.05148055978509729, 3, 1, 2.373057666666667, 7.119173, 0.03663411841135517, 0.03663411841135517, 8.789392685299333, 79.104534167694, 0.001885944923050103, 0.001885944923050103),
  @("M2", "Rgmb", "Neo1", "FAPs", 3, 1, 3.703826, 11.111478, 0.0514805597850973, 0.05148055978509729, 3, 1, 43.26393666666667, 129.79181, 0.6678877639810287, 0.6678877639810287, 160.2420934883533, 1442.17884139518, 0.0343832359633603, 0.0343832359633603),
  @("M2", "Rgmb", "Neo1", "M2", 3, 1, 3.703826, 11.111478, 0.0514805597850973, 0.05148055978509729, 3, 1, 0.4664283333333333, 1.399285, 0.007200495391983466, 0.007200495391983466, 1.727569388136666, 15.54812449323, 0.0003706855335093224, 0.0003706855335093224),
  @("M2", "Rgmb", "Neo1", "sCs", 3, 1, 3.703826, 11.111478, 0.0514805597850973, 0.05148055978509729, 3, 1, 18.67383333333333, 56.0215, 0.2882776222156329, 0.2882776222156329, 69.16462941966667, 622.481664777, 0.01484069336517758, 0.01484069336517758),
  @("sCs", "Rgmb", "Neo1", "ECs", 3, 1, 23.569813, 70.709439, 0.3276037176881591, 0.3276037176881591, 3, 1, 2.373057666666667, 7.119173, 0.03663411841135517, 0.03663411841135517, 55.93252544154966, 503.392728973947, 0.01200147338578819, 0.01200147338578819),
  @("sCs", "Rgmb", "Neo1", "FAPs", 3, 1, 23.569813, 70.709439, 0.3276037176881591, 0.3276037176881591, 3, 1, 43.26393666666667, 129.79181, 0.6678877639810287, 0.6678877639810287, 1019.722896877177, 9177.50607189459, 0.2188025144786168, 0.2188025144786168),
  @("sCs", "Rgmb", "Neo1", "M2", 3, 1, 23.569813, 70.709439, 0.3276037176881591, 0.3276037176881591, 3, 1, 0.4664283333333333, 1.399285, 0.007200495391983466, 0.007200495391983466, 10.99362859456833, 98.942657351115, 0.002358909059610242, 0.002358909059610242),
  @("sCs", "Rgmb", "Neo1", "sCs", 3, 1, 23.569813, 70.709439, 0.3276037176881591, 0.3276037176881591, 3, 1, 18.67383333333333, 56.0215, 0.2882776222156329, 0.2882776222156329, 440.1387596598333, 3961.248836938501, 0.09444082076414397, 0.09444082076414397)
)

$nRows = $data.Count
$nCols = $data[0].Count
$arr = New-Object 'object[,]' $nRows,$nCols
for ($i = 0; $i -lt $nRows; $i++) {
  for ($j = 0; $j -lt $nCols; $j++) {
    $arr[$i,$j] = $data[$i][$j]
  }
}

$ws.Range("A2:T17").Value = $arr
